# fixes for api & web demo tests
#
# Replaces the three GSMArena demo phones (Z2 / Galaxy On7 Pro / Galaxy S7
# active) with current Samsung Galaxy S10 family models, and leaves the
# Calculator sheet untouched. GSMArena becomes the active sheet/tab with
# H5 selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GSMArena")

# --- Row 2: GSMA0001 -------------------------------------------------
$ws1.Range("D2").ClearFormats()
$ws1.Range("D2").Value = "Galaxy S10+"
$ws1.Range("D2").Font.Name = "Menlo"
$ws1.Range("D2").Font.Color = 2236962
$ws1.Range("E2").Value = "6.4"""
$ws1.Range("F2").Value = "16MP"
$ws1.Range("G2").Value = "12GB RAM"
$ws1.Range("H2").Value = "4100mAh"

# --- Row 3: GSMA0002 -------------------------------------------------
$ws1.Range("D3").Value = "Galaxy S10"
$ws1.Range("E3").Value = "6.1"""
$ws1.Range("F3").Value = "16MP"
$ws1.Range("G3").Value = "8GB RAM"
$ws1.Range("H3").Value = "3400mAh"

# --- Row 4: GSMA0003 -------------------------------------------------
$ws1.Range("D4").Value = "Galaxy View2"
$ws1.Range("E4").Value = "17.3"""
$ws1.Range("F4").Value = "NO"
$ws1.Range("G4").Value = "3GB RAM"
$ws1.Range("H4").Value = "12000mAh"

# --- Make GSMArena the active/selected tab, cursor on H5 -------------
[void]$ws1.Activate()
[void]$ws1.Range("H5").Select()
